# Update the "last saved" datetime field placeholder text that appears on
# every slide's footer/date placeholder from 2019/2/21 -> 2019/2/23, then
# append two new slides (8 and 9) using the "Title and Content" layout,
# matching the author's "Add files via upload" commit.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Fix the cached datetimeFigureOut field text on every existing slide.
# ---------------------------------------------------------------------------
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if ($shape.HasTextFrame -and $shape.TextFrame.HasText) {
            $tr = $shape.TextFrame.TextRange
            if ($tr.Text -eq "2019/2/21") {
                $tr.Text = "2019/2/23"
            }
        }
    }
}

# ---------------------------------------------------------------------------
# 2) Append slide 8: "作成スケジュール"
# ---------------------------------------------------------------------------
$s8 = $p.Slides.Add($p.Slides.Count + 1, 16)

$title8 = $s8.Shapes.Item(1).TextFrame.TextRange
$title8.Text = "作成スケジュール"
$title8.Characters(1, 2).Text = "作成"

$body8 = $s8.Shapes.Item(2).TextFrame.TextRange
$body8.Text = "1.テーブル設計`r2.画面設計`r3.サーバサイド設計`r4.プログラミング`r…`r5.デプロイ方法"
$body8.ParagraphFormat.Bullet.Visible = 0

# ---------------------------------------------------------------------------
# 3) Append slide 9: "1.テーブル設計"
# ---------------------------------------------------------------------------
$s9 = $p.Slides.Add($p.Slides.Count + 1, 16)

$title9 = $s9.Shapes.Item(1).TextFrame.TextRange
$title9.Text = "1.テーブル設計"

$body9 = $s9.Shapes.Item(2).TextFrame.TextRange
$body9.Text = "SQLを使用する`r`rデータの更新、削除等を行うため。`r→それがなければNoSQLが良かった"
$body9.ParagraphFormat.Bullet.Visible = 0
